$wb = $excel.ActiveWorkbook

# The "Career Projection" sheet holds the header row that needs updating.
$ws = $wb.Worksheets.Item("Career Projection")

# Header row: "Promote" -> "New Rank", "Move" -> "New Zip"
# (set D1 first so the new shared strings are appended in the same
# order as the target workbook: "New Zip" then "New Rank")
$ws.Range("D1").Value = "New Zip"
$ws.Range("B1").Value = "New Rank"

# Make this sheet active/selected and move the selection to B2
$ws.Activate()
$ws.Range("B2").Select()
